$d = $word.ActiveDocument

# 1. "8 MIPS" -> "10 MIPS" (worst-case MIPS estimate correction)
$r1 = $d.Content
$r1.Find.Execute("Thus, 8 MIPS", $true, $false, $false, $false, $false, $true, 1, $false, "Thus, 10 MIPS", 1)

# 2. Remove the old _GoBack bookmark (it sat at the end of the "...CPU budget ... of
#    a single core." paragraph, right before the "Jinx scripts are designed..." paragraph).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 3. Rewrite the final sentence of the "Threaded Performance" section describing the
#    macOS / Dual Core behavior.
$r3 = $d.Content
$oldSentence = " and its Dual Core processor, per-thread performance actually drops once threads exceed the number of physical cores."
$newSentence = " and its Dual Core processor, per-thread performance drops significantly once the number of threads exceeds the number of physical cores."
$r3.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 1)

# 4. Re-insert the _GoBack bookmark right after "exceeds" in the newly rewritten sentence.
$r4 = $d.Content
$r4.Find.Execute("threads exceeds", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($r4.End, $r4.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Edits applied successfully"
